$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 47 - Year 2026
$ws.Range("B47").Value = 28.78
$ws.Range("C47").Value = 28.61
$ws.Range("D47").Value = 29.08
$ws.Range("E47").Value = 29.41
$ws.Range("F47").Value = 28.97
$ws.Range("G47").Value = 28.31
$ws.Range("H47").Value = 26.68
$ws.Range("I47").Value = 25.43
$ws.Range("J47").Value = 26.12
$ws.Range("K47").Value = 26.26
$ws.Range("L47").Value = 28.38
$ws.Range("M47").Value = 28.94

# Row 48 - Year 2027
$ws.Range("B48").Value = 28.4
$ws.Range("C48").Value = 28.71
$ws.Range("D48").Value = 28.88
$ws.Range("E48").Value = 29.35
$ws.Range("F48").Value = 29.32
$ws.Range("G48").Value = 27.93
$ws.Range("H48").Value = 27.41
$ws.Range("I48").Value = 25.37
$ws.Range("J48").Value = 26.29
$ws.Range("K48").Value = 27.08
$ws.Range("L48").Value = 27.77
$ws.Range("M48").Value = 28.33

# Row 49 - Year 2028
$ws.Range("B49").Value = 28.95
$ws.Range("C49").Value = 28.84
$ws.Range("D49").Value = 29.11
$ws.Range("E49").Value = 28.71
$ws.Range("F49").Value = 28.8
$ws.Range("G49").Value = 27.65
$ws.Range("H49").Value = 26.33
$ws.Range("I49").Value = 26.49
$ws.Range("J49").Value = 25.82
$ws.Range("K49").Value = 27.28
$ws.Range("L49").Value = 28.13
$ws.Range("M49").Value = 28.08
